$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data range before rewriting (old layout had data through E5)
$ws.Range("A1:E5").ClearContents()

# Write the transposed configuration table.
$ws.Range("A1").Value = "Department"
$ws.Range("B1").Value = "Mechanical Engineering"
$ws.Range("C1").Value = "Mechanical Engg"
$ws.Range("D1").Value = "Mech Eng"

$ws.Range("A2").Value = "Instructor"
$ws.Range("B2").Value = "Aryanci"
$ws.Range("C2").Value = "Cagri Aryanci"

$ws.Range("A3").Value = "Instructor"
$ws.Range("B3").Value = "Nobes"
$ws.Range("C3").Value = "David Nobes"
$ws.Range("D3").Value = "David S Nobes"
$ws.Range("E3").Value = "David S. Nobes"

$ws.Range("A4").Value = "Course Number"
$ws.Range("B4").Value = "MECE 260"
$ws.Range("C4").Value = "MecE 260"
$ws.Range("D4").Value = "MEC E 260"

$ws.Range("A5").Value = "Course Number"
$ws.Range("B5").Value = "MECE 265"
$ws.Range("C5").Value = "MecE 265"

# Column B needs to be widened to fit "Mechanical Engineering" / "Cagri Aryanci".
$ws.Columns.Item(2).ColumnWidth = 21.29

# Update the active selection to match the author's final cursor position.
$ws.Range("C5").Select()
